$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "autoawq" / "A" data row (original row 2), shifting the rest up.
$ws.Rows(2).Delete()

# Reflect the post-delete row selection (mirrors selecting the row after deleting it).
$ws.Range("A2:XFD2").Select()
